# Auto-generated edit script: updates currentAveragePrice / Leve price-profit
# columns (H-N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the refreshed
# market-price scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 173.66667
$ws.Range("J42").Value = 269.6
$ws.Range("L42").Value = 808.8000000000001
$ws.Range("N42").Value = -1268.8
# Row 106
$ws.Range("H106").Value = 1910.9333
$ws.Range("I106").Value = 1866.4615
$ws.Range("J106").Value = 2200
$ws.Range("K106").Value = 1866.4615
$ws.Range("L106").Value = 2200
$ws.Range("M106").Value = -1235.4615
$ws.Range("N106").Value = -3462
# Row 132
$ws.Range("H132").Value = 6640.7896
$ws.Range("I132").Value = 6954.1665
$ws.Range("K132").Value = 20862.4995
$ws.Range("M132").Value = -18332.4995
# Row 137
$ws.Range("H137").Value = 4264.722
$ws.Range("I137").Value = 2092.0588
$ws.Range("K137").Value = 6276.176399999999
$ws.Range("M137").Value = -3726.176399999999
# Row 138
$ws.Range("H138").Value = 5116.014
$ws.Range("J138").Value = 7673.3955
$ws.Range("L138").Value = 23020.1865
$ws.Range("N138").Value = -33300.1865

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 14287840
$ws.Range("I61").Value = 15626969
$ws.Range("J61").Value = 3790
$ws.Range("K61").Value = 15626969
$ws.Range("L61").Value = 3790
$ws.Range("M61").Value = -15626757
$ws.Range("N61").Value = -4214
# Row 63
$ws.Range("H63").Value = 1999
$ws.Range("J63").Value = 1999
$ws.Range("L63").Value = 1999
$ws.Range("N63").Value = -3371
# Row 66
$ws.Range("H66").Value = 1999
$ws.Range("J66").Value = 1999
$ws.Range("L66").Value = 9995
$ws.Range("N66").Value = -16859
# Row 132
$ws.Range("H132").Value = 25686182
$ws.Range("I132").Value = 11099.866
$ws.Range("J132").Value = 111269784
$ws.Range("K132").Value = 33299.598
$ws.Range("L132").Value = 333809352
$ws.Range("M132").Value = -30769.598
$ws.Range("N132").Value = -333814412
# Row 136
$ws.Range("H136").Value = 14287840
$ws.Range("I136").Value = 15626969
$ws.Range("J136").Value = 3790
$ws.Range("K136").Value = 46880907
$ws.Range("L136").Value = 11370
$ws.Range("M136").Value = -46878357
$ws.Range("N136").Value = -16470

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 9165
$ws.Range("I82").Value = 4874.1113
$ws.Range("K82").Value = 4874.1113
$ws.Range("M82").Value = -4491.1113
# Row 85
$ws.Range("H85").Value = 9165
$ws.Range("I85").Value = 4874.1113
$ws.Range("K85").Value = 4874.1113
$ws.Range("M85").Value = -3548.1113
# Row 99
$ws.Range("H99").Value = 5122.3335
$ws.Range("J99").Value = 6794.8335
$ws.Range("L99").Value = 6794.8335
$ws.Range("N99").Value = -9790.833500000001
# Row 105
$ws.Range("H105").Value = 7736.615
$ws.Range("I105").Value = 15739.875
$ws.Range("K105").Value = 15739.875
$ws.Range("M105").Value = -13992.875
# Row 134
$ws.Range("H134").Value = 3393.8667
$ws.Range("I134").Value = 3488.162
$ws.Range("J134").Value = 2957.75
$ws.Range("K134").Value = 10464.486
$ws.Range("L134").Value = 8873.25
$ws.Range("M134").Value = -7929.485999999999
$ws.Range("N134").Value = -13943.25

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 190.08333
$ws.Range("I19").Value = 112.14286
$ws.Range("K19").Value = 112.14286
$ws.Range("M19").Value = 57.85714
# Row 24
$ws.Range("H24").Value = 190.08333
$ws.Range("I24").Value = 112.14286
$ws.Range("K24").Value = 112.14286
$ws.Range("M24").Value = 57.85714
# Row 31
$ws.Range("H31").Value = 10004692
$ws.Range("I31").Value = 3074.125
$ws.Range("K31").Value = 3074.125
$ws.Range("M31").Value = -2779.125
# Row 34
$ws.Range("H34").Value = 10004692
$ws.Range("I34").Value = 3074.125
$ws.Range("K34").Value = 3074.125
$ws.Range("M34").Value = -2872.125
# Row 58
$ws.Range("H58").Value = 2379.3333
$ws.Range("I58").Value = 2409
$ws.Range("J58").Value = 2320
$ws.Range("K58").Value = 2409
$ws.Range("L58").Value = 2320
$ws.Range("M58").Value = -2206
$ws.Range("N58").Value = -2726
# Row 62
$ws.Range("H62").Value = 3407.9167
$ws.Range("I62").Value = 3372.2727
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 3372.2727
$ws.Range("L62").Value = 3800
$ws.Range("M62").Value = -2748.2727
$ws.Range("N62").Value = -5048
# Row 65
$ws.Range("H65").Value = 3407.9167
$ws.Range("I65").Value = 3372.2727
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 16861.3635
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = -13741.3635
$ws.Range("N65").Value = -25240
# Row 136
$ws.Range("H136").Value = 2379.3333
$ws.Range("I136").Value = 2409
$ws.Range("J136").Value = 2320
$ws.Range("K136").Value = 7227
$ws.Range("L136").Value = 6960
$ws.Range("M136").Value = -4677
$ws.Range("N136").Value = -12060

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 497.5
$ws.Range("I68").Value = 497.5
$ws.Range("K68").Value = 1492.5
$ws.Range("M68").Value = -681.5
# Row 71
$ws.Range("H71").Value = 497.5
$ws.Range("I71").Value = 497.5
$ws.Range("K71").Value = 4477.5
$ws.Range("M71").Value = -421.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 45
$ws.Range("H45").Value = 40999.75
$ws.Range("I45").Value = 39999.5
$ws.Range("K45").Value = 39999.5
$ws.Range("M45").Value = -39440.5
# Row 51
$ws.Range("H51").Value = 79999.5
$ws.Range("I51").Value = 79999.5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 79999.5
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -79490.5
$ws.Range("N51").Value = $null
# Row 126
$ws.Range("H126").Value = 204493.8
$ws.Range("I126").Value = 500750
$ws.Range("J126").Value = 6989.6665
$ws.Range("K126").Value = 1502250
$ws.Range("L126").Value = 20968.9995
$ws.Range("M126").Value = -1499780
$ws.Range("N126").Value = -25908.9995
# Row 131
$ws.Range("H131").Value = 99999.5
$ws.Range("J131").Value = 99999.5
$ws.Range("L131").Value = 99999.5
$ws.Range("N131").Value = -110079.5
# Row 132
$ws.Range("H132").Value = 5104.0435
$ws.Range("I132").Value = 5104.4546
$ws.Range("J132").Value = 5095
$ws.Range("K132").Value = 15313.3638
$ws.Range("L132").Value = 15285
$ws.Range("M132").Value = -12783.3638
$ws.Range("N132").Value = -20345

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 5000
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5226
# Row 22
$ws.Range("H22").Value = 2225.5881
$ws.Range("I22").Value = 2205
$ws.Range("J22").Value = 2255
$ws.Range("K22").Value = 2205
$ws.Range("L22").Value = 2255
$ws.Range("M22").Value = -1910
$ws.Range("N22").Value = -2845
# Row 25
$ws.Range("H25").Value = 3999.3333
$ws.Range("I25").Value = 3990
$ws.Range("K25").Value = 3990
$ws.Range("M25").Value = -3760
# Row 27
$ws.Range("H27").Value = 2225.5881
$ws.Range("I27").Value = 2205
$ws.Range("J27").Value = 2255
$ws.Range("K27").Value = 2205
$ws.Range("L27").Value = 2255
$ws.Range("M27").Value = -2098
$ws.Range("N27").Value = -2469
# Row 28
$ws.Range("H28").Value = 5000
$ws.Range("J28").Value = 5000
$ws.Range("L28").Value = 5000
$ws.Range("N28").Value = -5464
# Row 37
$ws.Range("H37").Value = 5000
$ws.Range("J37").Value = 5000
$ws.Range("L37").Value = 5000
$ws.Range("N37").Value = -5214
# Row 46
$ws.Range("H46").Value = 1454.8572
$ws.Range("I46").Value = 1030.6666
$ws.Range("K46").Value = 1030.6666
$ws.Range("M46").Value = -842.6666
# Row 55
$ws.Range("H55").Value = 564.5599999999999
$ws.Range("I55").Value = 321.05264
$ws.Range("J55").Value = 1335.6666
$ws.Range("K55").Value = 321.05264
$ws.Range("L55").Value = 1335.6666
$ws.Range("M55").Value = -148.05264
$ws.Range("N55").Value = -1681.6666
# Row 136
$ws.Range("H136").Value = 1543484.2
$ws.Range("J136").Value = 15750
$ws.Range("L136").Value = 47250
$ws.Range("N136").Value = -52350

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 77693570
$ws.Range("I100").Value = 112223510
$ws.Range("K100").Value = 224447020
$ws.Range("M100").Value = -224446479
# Row 107
$ws.Range("H107").Value = 2485.3333
$ws.Range("I107").Value = 968.5
$ws.Range("J107").Value = 3243.75
$ws.Range("K107").Value = 2905.5
$ws.Range("L107").Value = 9731.25
$ws.Range("M107").Value = -985.5
$ws.Range("N107").Value = -13571.25
# Row 132
$ws.Range("H132").Value = 3557.75
$ws.Range("I132").Value = 3299.3333
$ws.Range("K132").Value = 9897.999899999999
$ws.Range("M132").Value = -7367.999899999999
